$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data runs through row 269 (date serial 44343, 2021-05-27).
# Append daily rows through row 301 (date serial 44375, 2021-06-28),
# matching the existing pattern: col A = date serial, cols B/C/D = 0.
$startRow = 270
$endRow = 301
$startSerial = 44344

for ($row = $startRow; $row -le $endRow; $row++) {
    $serial = $startSerial + ($row - $startRow)
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

# Replicate the date-column formatting/style used by the preceding rows
# (column A uses style index 2 in the original file: centered, bordered,
# bold, date-time number format).
$ws.Range("A269").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)
